# Atualização automática: 2025-08-08 10:00:35
# Appends one new data row (row 18) to Sheet1, mirroring the existing
# "Fly_ID" tracking table rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18

$ws.Cells.Item($row, 1).Value = "ce49ea82-267b-4cc5-86c9-7e4337c56079"
$ws.Cells.Item($row, 2).Value = "mosca"

# Column C holds a date serial, styled the same way as the row above it
# (custom date/time number format already used throughout the column).
$ws.Cells.Item($row, 3).Value = 45877
$ws.Cells.Item($row, 3).NumberFormat = $ws.Cells.Item($row - 1, 3).NumberFormat

$ws.Cells.Item($row, 4).Value = "image_20250808100711_ppp0.jpg"
$ws.Cells.Item($row, 5).Value = "PLACA_20250717165933"
$ws.Cells.Item($row, 6).Value = "Beja"
$ws.Cells.Item($row, 7).Value = 38.02035
$ws.Cells.Item($row, 8).Value = -7.94715
$ws.Cells.Item($row, 9).Value = "1182,409,1232,451"

# Column J stores confidence values as plain text (e.g. "0.75", "0.70")
# everywhere else in the sheet, so force text with a leading apostrophe
# rather than letting it be auto-detected as a number.
$ws.Cells.Item($row, 10).Value = "'0.75"
